$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 10 (Objetivos): only the B/C value text changes; label + style
# stay the same, so just overwrite the cell values in place.
# ------------------------------------------------------------------
$ws.Range("B10:C10").Value = "Introduzir os conceitos fundamentais da ciência administração e de configurações de uma organização."

# ------------------------------------------------------------------
# Rows 12-21 get fully restructured: a row is inserted right after
# "Docentes responsaveis:" to hold its value, and a new row is
# appended at the end for "Bibliografia:". Clear the whole block
# first (dropping stale cells/heights), then rebuild rows 12-22 from
# scratch, copying cell formatting from the untouched template row 3
# (A3 = label style, B3 = plain value style, C3 = red value style).
# ------------------------------------------------------------------
$ws.Range("A12:C21").EntireRow.Delete()

function Set-LabelCell($addr, $text) {
    $ws.Range("A3").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $ws.Range($addr).Value = $text
}
function Set-PlainValueCell($addr, $text) {
    $ws.Range("B3").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $ws.Range($addr).Value = $text
}
function Set-RedValueCell($addr, $text) {
    $ws.Range("C3").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $ws.Range($addr).Value = $text
}

# --- Row 12 ---
Set-LabelCell "A12" "Docentes responsáveis:"

# --- Row 13 ---
Set-PlainValueCell "B13" "11079086 - Herlandí de Souza Andrade"
Set-RedValueCell "C13" "11079086 - Herlandí de Souza Andrade"

# --- Row 14 ---
Set-LabelCell "A14" "Programa resumido:"
Set-PlainValueCell "B14" "1. Áreas de Atuação da Administração.2. Estrutura organizacional."
Set-RedValueCell "C14" "1. Áreas de Atuação da Administração.2. Estrutura organizacional."
$ws.Rows(14).RowHeight = 60

# --- Row 15 ---
Set-LabelCell "A15" "Short syllabus:"
Set-PlainValueCell "B15" "1. Management Practice Areas. 2. Organizational structure"
Set-RedValueCell "C15" "1. Management Practice Areas. 2. Organizational structure"
$ws.Rows(15).RowHeight = 60

# --- Row 16 ---
Set-LabelCell "A16" "Programa:"
Set-PlainValueCell "B16" "1. Noções básicas de Marketing, Finanças e Recursos Humanos. 2. Diferentes configurações de organização."
Set-RedValueCell "C16" "1. Noções básicas de Marketing, Finanças e Recursos Humanos. 2. Diferentes configurações de organização."
$ws.Rows(16).RowHeight = 120

# --- Row 17 ---
Set-LabelCell "A17" "Syllabus:"
Set-PlainValueCell "B17" "1. Basic notions of Marketing, Finance and Human Resources.2. Different organization settings."
Set-RedValueCell "C17" "1. Basic notions of Marketing, Finance and Human Resources.2. Different organization settings."
$ws.Rows(17).RowHeight = 120

# --- Row 18 ---
Set-LabelCell "A18" "Avaliação:"

# --- Row 19 ---
Set-LabelCell "A19" "Método:"
Set-PlainValueCell "B19" "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."
Set-RedValueCell "C19" "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."
$ws.Rows(19).RowHeight = 60

# --- Row 20 ---
Set-LabelCell "A20" "Critério:"
Set-PlainValueCell "B20" "Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas"
Set-RedValueCell "C20" "Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas"
$ws.Rows(20).RowHeight = 60

# --- Row 21 ---
Set-LabelCell "A21" "Norma de recuperação:"
Set-PlainValueCell "B21" "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação"
Set-RedValueCell "C21" "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação"
$ws.Rows(21).RowHeight = 60

# --- Row 22 ---
Set-LabelCell "A22" "Bibliografia:"
Set-PlainValueCell "B22" "Chiavenato, I. Gestão de Pessoas. 4 ed. São Paulo: Manole, 2014.Chiavenato, I. Recursos Humanos: o capital humano das organizações. 10 ed. Rio de Janeiro, Campus, 2015.ROBBINS, S. P.; DECENZO, D. A.; WOLTER, R. Fundamentos de Gestão de Pessoas. São Paulo, saraiva, 2013.KOTLER, P. - ARMSTRONG, G. Princípios De Marketing. 15 ed. São Paulo: Pearson, 2014.KOTLER, P.; KELLER, K. L. Administração De Marketing. 15 ed. São Paulo: Pearson, 2019.CHIAVENATO, I. Introdução À Teoria Geral da Administração. 9 ed. São Paulo: Manole, 2014. MAXIMIANO, A. C. A. Teoria Geral da Administração: da Revolução Urbana À Revolução Digital. 8 ed. São Paulo: Atlas, 2017.GUERRINI, F. M.; ESCRIÇÃO FILHO, E.; ROSIM, D. Administração Para Engenheiros. Rio de Janeiro: Campus, 2016.CHIAVENATO, I. Administração Para Não Administradores: a Gestão de Negócios Ao Alcance de Todos. 2 ed. São Paulo: Manole, 2011.SILVA, M. M. L. Administração para Estudantes e Profissionais de Áreas Técnicas. São Paulo: Brasport, 2018.GITMAN, L. J. - ZUTTER, C. J. Princípios de Administração Financeira. 14 ed. São Paulo: Perason, 2017.GROPPELLI, A. A.; NIKBAKHT, E. Administração Financeira. 3 ed. São Paulo: Saraiva, 2010.MARCOUSÉ, I.; SURRIDGE, M.; GILLESPIE, A. Finanças. São Paulo: Saraiva, 2013.BOLMAN, L.G.; DEAL, T.E. Reframing organizations. San Francisco, John Wiley, 2013KOTLER, P.. O Marketing sem segredos. 1 ed. Porto Alegre. Bookman, 2005MINTZBERG, H. Criando organizações eficazes. 2 ed. São Paulo, Atlas, 2006.MORGAN, G. Imagens da organização. São Paulo, Atlas, 1996."
Set-RedValueCell "C22" "Chiavenato, I. Gestão de Pessoas. 4 ed. São Paulo: Manole, 2014.Chiavenato, I. Recursos Humanos: o capital humano das organizações. 10 ed. Rio de Janeiro, Campus, 2015.ROBBINS, S. P.; DECENZO, D. A.; WOLTER, R. Fundamentos de Gestão de Pessoas. São Paulo, saraiva, 2013.KOTLER, P. - ARMSTRONG, G. Princípios De Marketing. 15 ed. São Paulo: Pearson, 2014.KOTLER, P.; KELLER, K. L. Administração De Marketing. 15 ed. São Paulo: Pearson, 2019.CHIAVENATO, I. Introdução À Teoria Geral da Administração. 9 ed. São Paulo: Manole, 2014. MAXIMIANO, A. C. A. Teoria Geral da Administração: da Revolução Urbana À Revolução Digital. 8 ed. São Paulo: Atlas, 2017.GUERRINI, F. M.; ESCRIÇÃO FILHO, E.; ROSIM, D. Administração Para Engenheiros. Rio de Janeiro: Campus, 2016.CHIAVENATO, I. Administração Para Não Administradores: a Gestão de Negócios Ao Alcance de Todos. 2 ed. São Paulo: Manole, 2011.SILVA, M. M. L. Administração para Estudantes e Profissionais de Áreas Técnicas. São Paulo: Brasport, 2018.GITMAN, L. J. - ZUTTER, C. J. Princípios de Administração Financeira. 14 ed. São Paulo: Perason, 2017.GROPPELLI, A. A.; NIKBAKHT, E. Administração Financeira. 3 ed. São Paulo: Saraiva, 2010.MARCOUSÉ, I.; SURRIDGE, M.; GILLESPIE, A. Finanças. São Paulo: Saraiva, 2013.BOLMAN, L.G.; DEAL, T.E. Reframing organizations. San Francisco, John Wiley, 2013KOTLER, P.. O Marketing sem segredos. 1 ed. Porto Alegre. Bookman, 2005MINTZBERG, H. Criando organizações eficazes. 2 ed. São Paulo, Atlas, 2006.MORGAN, G. Imagens da organização. São Paulo, Atlas, 1996."
$ws.Rows(22).RowHeight = 120

$excel.CutCopyMode = $false